$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

# Remove existing hyperlinks so we can rebuild them cleanly against the new row layout
$ws.Cells.Hyperlinks.Delete()

# Row 1 (header) is unchanged.

# Row 2
$ws.Range("A2").Value = '2025-10-17 01:16:27'
$ws.Range("B2").Value = '【Azure/RAG】社内文書検索AIチャットボットの精度向上&内製化支援パートナー募集!'
$ws.Range("C2").Value = 'システム開発'
$ws.Range("D2").Value = '300,000 円 ~ 500,000 円 / 固定'
$ws.Range("E2").Value = '期限情報なし'
$ws.Range("F2").Value = 'https://www.lancers.jp/work/detail/5413954'
$ws.Range("G2").Value = 310
$ws.Range("H2").Value = '🔥AI,Ai'

# Row 3
$ws.Range("A3").Value = '2025-10-17 01:16:27'
$ws.Range("B3").Value = '生成AIの技術顧問を募集!事業の技術選定をリードするAI専門家を募集! 【週1日〜/フルリモート】'
$ws.Range("C3").Value = 'システム開発'
$ws.Range("D3").Value = '200,000 円 ~ 300,000 円 / 固定'
$ws.Range("E3").Value = '期限情報なし'
$ws.Range("F3").Value = 'https://www.lancers.jp/work/detail/5413955'
$ws.Range("G3").Value = 303
$ws.Range("H3").Value = '🔥AI,Ai'

# Row 4
$ws.Range("A4").Value = '2025-10-17 01:16:27'
$ws.Range("B4").Value = '【せどり×ツール製作】APIを使用したせどりツールを製作できるエンジニアさんを募集します♪'
$ws.Range("C4").Value = 'システム開発'
$ws.Range("D4").Value = '20,000 円 ~ 50,000 円 / 固定'
$ws.Range("E4").Value = '期限情報なし'
$ws.Range("F4").Value = 'https://www.lancers.jp/work/detail/5217096'
$ws.Range("G4").Value = 243
$ws.Range("H4").Value = '🔥API ◆ツール'

# Row 5
$ws.Range("A5").Value = '2025-10-17 01:16:27'
$ws.Range("B5").Value = '【募集】RPAツール「RoboTANGO」設定代行の専門家を探しています'
$ws.Range("C5").Value = 'システム開発'
$ws.Range("D5").Value = '50,000 円 ~ 100,000 円 / 固定'
$ws.Range("E5").Value = '期限情報なし'
$ws.Range("F5").Value = 'https://www.lancers.jp/work/detail/5405023'
$ws.Range("G5").Value = 178
$ws.Range("H5").Value = '★bot ◆ツール'

# Row 6
$ws.Range("A6").Value = '2025-10-17 01:16:27'
$ws.Range("B6").Value = '【急募】キントーン見積をExcelに変換するツール開発'
$ws.Range("C6").Value = 'システム開発'
$ws.Range("D6").Value = '10,000 円 ~ 20,000 円 / 固定'
$ws.Range("E6").Value = '期限情報なし'
$ws.Range("F6").Value = 'https://www.lancers.jp/work/detail/5414167'
$ws.Range("G6").Value = 120
$ws.Range("H6").Value = '◆ツール,開発'

# Row 7
$ws.Range("A7").Value = '2025-10-17 01:16:27'
$ws.Range("B7").Value = '大手クレジットカード企業向け、Google Cloudを利用したアジャイル開発共通基盤案件_ワーカー'
$ws.Range("C7").Value = 'システム開発'
$ws.Range("D7").Value = '500,000 円 ~ 1,000,000 円 / 固定'
$ws.Range("E7").Value = '期限情報なし'
$ws.Range("F7").Value = 'https://www.lancers.jp/work/detail/5414354'
$ws.Range("G7").Value = 75
$ws.Range("H7").Value = '◆開発'

# Row 8
$ws.Range("A8").Value = '2025-10-17 01:16:27'
$ws.Range("B8").Value = '大手クレジットカード企業向け、Google Cloudを利用したアジャイル開発共通基盤案件'
$ws.Range("C8").Value = 'システム開発'
$ws.Range("D8").Value = '500,000 円 ~ 1,000,000 円 / 固定'
$ws.Range("E8").Value = '期限情報なし'
$ws.Range("F8").Value = 'https://www.lancers.jp/work/detail/5414353'
$ws.Range("G8").Value = 75
$ws.Range("H8").Value = '◆開発'

# Row 9
$ws.Range("A9").Value = '2025-10-17 01:16:27'
$ws.Range("B9").Value = '【長期依頼】海外クリエイター向けサービスの開発保守|Laravel+Livewireエンジニア募集'
$ws.Range("C9").Value = 'システム開発'
$ws.Range("D9").Value = '300,000 円 ~ 500,000 円 / 固定'
$ws.Range("E9").Value = '期限情報なし'
$ws.Range("F9").Value = 'https://www.lancers.jp/work/detail/5414105'
$ws.Range("G9").Value = 75
$ws.Range("H9").Value = '◆開発'

# Row 10
$ws.Range("A10").Value = '2025-10-17 01:16:27'
$ws.Range("B10").Value = '【急募】見積書自動作成機能の開発'
$ws.Range("C10").Value = 'システム開発'
$ws.Range("D10").Value = '500,000 円 ~ 1,000,000 円 / 固定'
$ws.Range("E10").Value = '期限情報なし'
$ws.Range("F10").Value = 'https://www.lancers.jp/work/detail/5414108'
$ws.Range("G10").Value = 75
$ws.Range("H10").Value = '◆開発'

# Row 11
$ws.Range("A11").Value = '2025-10-17 01:16:27'
$ws.Range("B11").Value = 'セレニウムを用いた自動発注ツールの修正・機能追加'
$ws.Range("C11").Value = 'システム開発'
$ws.Range("D11").Value = '50,000 円 ~ 100,000 円 / 固定'
$ws.Range("E11").Value = '期限情報なし'
$ws.Range("F11").Value = 'https://www.lancers.jp/work/detail/5413916'
$ws.Range("G11").Value = 73
$ws.Range("H11").Value = '◆ツール'

# Row 12
$ws.Range("A12").Value = '2025-10-17 01:16:27'
$ws.Range("B12").Value = 'Symfoware/RDB向け参照系SQLクエリ実行ライブラリの開発依頼'
$ws.Range("C12").Value = 'システム開発'
$ws.Range("D12").Value = '50,000 円 ~ 100,000 円 / 固定'
$ws.Range("E12").Value = '期限情報なし'
$ws.Range("F12").Value = 'https://www.lancers.jp/work/detail/5414368'
$ws.Range("G12").Value = 68
$ws.Range("H12").Value = '◆開発'

# Row 13
$ws.Range("A13").Value = '2025-10-17 01:16:27'
$ws.Range("B13").Value = '【音声コマンド起動】超小型・低電力レコーダーのプロトタイプ開発'
$ws.Range("C13").Value = 'システム開発'
$ws.Range("D13").Value = '20,000 円 ~ 50,000 円 / 固定'
$ws.Range("E13").Value = '期限情報なし'
$ws.Range("F13").Value = 'https://www.lancers.jp/work/detail/5413958'
$ws.Range("G13").Value = 63
$ws.Range("H13").Value = '◆開発'

# Row 14
$ws.Range("A14").Value = '2025-10-17 01:16:27'
$ws.Range("B14").Value = 'リアルタイム音声チャットボット強化・管理UI構築(ASR/LLM/TTS最適化対応)'
$ws.Range("C14").Value = 'システム開発'
$ws.Range("D14").Value = '1,000,000 円 ~ 3,000,000 円 / 固定'
$ws.Range("E14").Value = '期限情報なし'
$ws.Range("F14").Value = 'https://www.lancers.jp/work/detail/5414569'
$ws.Range("G14").Value = 45
$ws.Range("H14").Value = '◇管理'

# Row 15
$ws.Range("A15").Value = '2025-10-17 01:16:27'
$ws.Range("B15").Value = 'PowerAutomateメール監視して件名と本文内の条件一致時、社内システム操作&メール転送したい'
$ws.Range("C15").Value = 'システム開発'
$ws.Range("D15").Value = '20,000 円 ~ 50,000 円 / 固定'
$ws.Range("E15").Value = '期限情報なし'
$ws.Range("F15").Value = 'https://www.lancers.jp/work/detail/5414579'
$ws.Range("G15").Value = 28

# Row 16
$ws.Range("A16").Value = '2025-10-17 01:16:27'
$ws.Range("B16").Value = '初回 【急募・即決します】VBAで1問1答問題集の作成'
$ws.Range("C16").Value = 'システム開発'
$ws.Range("D16").Value = '~ 5,000 円 / 固定'
$ws.Range("E16").Value = '期限情報なし'
$ws.Range("F16").Value = 'https://www.lancers.jp/work/detail/5414812'
$ws.Range("G16").Value = 10

# Re-add hyperlinks for the URL column (F2:F16), matching each row's target cell
$ws.Hyperlinks.Add($ws.Range("F2"), 'https://www.lancers.jp/work/detail/5413954') | Out-Null
$ws.Range("F2").Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Range("F3"), 'https://www.lancers.jp/work/detail/5413955') | Out-Null
$ws.Range("F3").Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Range("F4"), 'https://www.lancers.jp/work/detail/5217096') | Out-Null
$ws.Range("F4").Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Range("F5"), 'https://www.lancers.jp/work/detail/5405023') | Out-Null
$ws.Range("F5").Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Range("F6"), 'https://www.lancers.jp/work/detail/5414167') | Out-Null
$ws.Range("F6").Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Range("F7"), 'https://www.lancers.jp/work/detail/5414354') | Out-Null
$ws.Range("F7").Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Range("F8"), 'https://www.lancers.jp/work/detail/5414353') | Out-Null
$ws.Range("F8").Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Range("F9"), 'https://www.lancers.jp/work/detail/5414105') | Out-Null
$ws.Range("F9").Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Range("F10"), 'https://www.lancers.jp/work/detail/5414108') | Out-Null
$ws.Range("F10").Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Range("F11"), 'https://www.lancers.jp/work/detail/5413916') | Out-Null
$ws.Range("F11").Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Range("F12"), 'https://www.lancers.jp/work/detail/5414368') | Out-Null
$ws.Range("F12").Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Range("F13"), 'https://www.lancers.jp/work/detail/5413958') | Out-Null
$ws.Range("F13").Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Range("F14"), 'https://www.lancers.jp/work/detail/5414569') | Out-Null
$ws.Range("F14").Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Range("F15"), 'https://www.lancers.jp/work/detail/5414579') | Out-Null
$ws.Range("F15").Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Range("F16"), 'https://www.lancers.jp/work/detail/5414812') | Out-Null
$ws.Range("F16").Style = "Hyperlink"

Write-Output "done"
